# Update NATMI ligand-receptor TPM values (Tgfb1-Eng) with recomputed
# expression / specificity / edge-weight figures for rows 2-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 54.53585066666667
$ws.Range("H2").Value = 163.607552
$ws.Range("I2").Value = 0.3031388658437607
$ws.Range("J2").Value = 0.3031388658437607
$ws.Range("M2").Value = 247.0944516666667
$ws.Range("N2").Value = 741.283355
$ws.Range("O2").Value = 0.8050739182622993
$ws.Range("P2").Value = 0.8050739182622993
$ws.Range("Q2").Value = 13475.50611665522
$ws.Range("R2").Value = 121279.555049897
$ws.Range("S2").Value = 0.2440491945024259
$ws.Range("T2").Value = 0.2440491945024259

$ws.Range("G3").Value = 54.53585066666667
$ws.Range("H3").Value = 163.607552
$ws.Range("I3").Value = 0.3031388658437607
$ws.Range("J3").Value = 0.3031388658437607
$ws.Range("O3").Value = 0.1379009747488701
$ws.Range("P3").Value = 0.13790097474887
$ws.Range("Q3").Value = 2308.217154434845
$ws.Range("R3").Value = 20773.9543899136
$ws.Range("S3").Value = 0.04180314508412154
$ws.Range("T3").Value = 0.04180314508412154

$ws.Range("G4").Value = 54.53585066666667
$ws.Range("H4").Value = 163.607552
$ws.Range("I4").Value = 0.3031388658437607
$ws.Range("J4").Value = 0.3031388658437607
$ws.Range("M4").Value = 11.590146
$ws.Range("N4").Value = 34.770438
$ws.Range("O4").Value = 0.03776258103132013
$ws.Range("P4").Value = 0.03776258103132013
$ws.Range("Q4").Value = 632.078471460864
$ws.Range("R4").Value = 5688.706243147776
$ws.Range("S4").Value = 0.01144730598516749
$ws.Range("T4").Value = 0.0114473059851675

$ws.Range("G5").Value = 54.53585066666667
$ws.Range("H5").Value = 163.607552
$ws.Range("I5").Value = 0.3031388658437607
$ws.Range("J5").Value = 0.3031388658437607
$ws.Range("M5").Value = 5.912082333333333
$ws.Range("N5").Value = 17.736247
$ws.Range("O5").Value = 0.01926252595751047
$ws.Range("P5").Value = 0.01926252595751047
$ws.Range("Q5").Value = 322.4204392597049
$ws.Range("R5").Value = 2901.783953337344
$ws.Range("S5").Value = 0.005839220272045724
$ws.Range("T5").Value = 0.005839220272045725

$ws.Range("I6").Value = 0.1026363515063155
$ws.Range("J6").Value = 0.1026363515063155
$ws.Range("M6").Value = 247.0944516666667
$ws.Range("N6").Value = 741.283355
$ws.Range("O6").Value = 0.8050739182622993
$ws.Range("P6").Value = 0.8050739182622993
$ws.Range("Q6").Value = 4562.518826693025
$ws.Range("R6").Value = 41062.66944023723
$ws.Range("S6").Value = 0.08262984966333606
$ws.Range("T6").Value = 0.08262984966333607

$ws.Range("I7").Value = 0.1026363515063155
$ws.Range("J7").Value = 0.1026363515063155
$ws.Range("O7").Value = 0.1379009747488701
$ws.Range("P7").Value = 0.13790097474887
$ws.Range("S7").Value = 0.01415365291738856
$ws.Range("T7").Value = 0.01415365291738856

$ws.Range("I8").Value = 0.1026363515063155
$ws.Range("J8").Value = 0.1026363515063155
$ws.Range("M8").Value = 11.590146
$ws.Range("N8").Value = 34.770438
$ws.Range("O8").Value = 0.03776258103132013
$ws.Range("P8").Value = 0.03776258103132013
$ws.Range("Q8").Value = 214.008282955932
$ws.Range("R8").Value = 1926.074546603388
$ws.Range("S8").Value = 0.003875813540516295
$ws.Range("T8").Value = 0.003875813540516295

$ws.Range("I9").Value = 0.1026363515063155
$ws.Range("J9").Value = 0.1026363515063155
$ws.Range("M9").Value = 5.912082333333333
$ws.Range("N9").Value = 17.736247
$ws.Range("O9").Value = 0.01926252595751047
$ws.Range("P9").Value = 0.01926252595751047
$ws.Range("Q9").Value = 109.1646808289358
$ws.Range("R9").Value = 982.4821274604219
$ws.Range("S9").Value = 0.001977035385074571
$ws.Range("T9").Value = 0.001977035385074571

$ws.Range("G10").Value = 12.55635966666667
$ws.Range("H10").Value = 37.669079
$ws.Range("I10").Value = 0.06979483370938171
$ws.Range("J10").Value = 0.06979483370938172
$ws.Range("M10").Value = 247.0944516666667
$ws.Range("N10").Value = 741.283355
$ws.Range("O10").Value = 0.8050739182622993
$ws.Range("P10").Value = 0.8050739182622993
$ws.Range("Q10").Value = 3102.606806764449
$ws.Range("R10").Value = 27923.46126088004
$ws.Range("S10").Value = 0.05619000024887754
$ws.Range("T10").Value = 0.05619000024887755

$ws.Range("G11").Value = 12.55635966666667
$ws.Range("H11").Value = 37.669079
$ws.Range("I11").Value = 0.06979483370938171
$ws.Range("J11").Value = 0.06979483370938172
$ws.Range("O11").Value = 0.1379009747488701
$ws.Range("P11").Value = 0.13790097474887
$ws.Range("Q11").Value = 531.4449930744111
$ws.Range("R11").Value = 4783.0049376697
$ws.Range("S11").Value = 0.009624775600959031
$ws.Range("T11").Value = 0.009624775600959031

$ws.Range("G12").Value = 12.55635966666667
$ws.Range("H12").Value = 37.669079
$ws.Range("I12").Value = 0.06979483370938171
$ws.Range("J12").Value = 0.06979483370938172
$ws.Range("M12").Value = 11.590146
$ws.Range("N12").Value = 34.770438
$ws.Range("O12").Value = 0.03776258103132013
$ws.Range("P12").Value = 0.03776258103132013
$ws.Range("Q12").Value = 145.530041765178
$ws.Range("R12").Value = 1309.770375886602
$ws.Range("S12").Value = 0.00263563306351804
$ws.Range("T12").Value = 0.002635633063518041

$ws.Range("G13").Value = 12.55635966666667
$ws.Range("H13").Value = 37.669079
$ws.Range("I13").Value = 0.06979483370938171
$ws.Range("J13").Value = 0.06979483370938172
$ws.Range("M13").Value = 5.912082333333333
$ws.Range("N13").Value = 17.736247
$ws.Range("O13").Value = 0.01926252595751047
$ws.Range("P13").Value = 0.01926252595751047
$ws.Range("Q13").Value = 74.23423215627922
$ws.Range("R13").Value = 668.1080894065129
$ws.Range("S13").Value = 0.001344424796027092
$ws.Range("T13").Value = 0.001344424796027092

$ws.Range("G14").Value = 94.34696966666667
$ws.Range("H14").Value = 283.040909
$ws.Range("I14").Value = 0.524429948940542
$ws.Range("J14").Value = 0.5244299489405421
$ws.Range("M14").Value = 247.0944516666667
$ws.Range("N14").Value = 741.283355
$ws.Range("O14").Value = 0.8050739182622993
$ws.Range("P14").Value = 0.8050739182622993
$ws.Range("Q14").Value = 23312.61273619664
$ws.Range("R14").Value = 209813.5146257697
$ws.Range("S14").Value = 0.4222048738476597
$ws.Range("T14").Value = 0.4222048738476598

$ws.Range("G15").Value = 94.34696966666667
$ws.Range("H15").Value = 283.040909
$ws.Range("I15").Value = 0.524429948940542
$ws.Range("J15").Value = 0.5244299489405421
$ws.Range("O15").Value = 0.1379009747488701
$ws.Range("P15").Value = 0.13790097474887
$ws.Range("Q15").Value = 3993.213476848744
$ws.Range("R15").Value = 35938.9212916387
$ws.Range("S15").Value = 0.0723194011464009
$ws.Range("T15").Value = 0.0723194011464009

$ws.Range("G16").Value = 94.34696966666667
$ws.Range("H16").Value = 283.040909
$ws.Range("I16").Value = 0.524429948940542
$ws.Range("J16").Value = 0.5244299489405421
$ws.Range("M16").Value = 11.590146
$ws.Range("N16").Value = 34.770438
$ws.Range("O16").Value = 0.03776258103132013
$ws.Range("P16").Value = 0.03776258103132013
$ws.Range("Q16").Value = 1093.495153094238
$ws.Range("R16").Value = 9841.456377848142
$ws.Range("S16").Value = 0.0198038284421183
$ws.Range("T16").Value = 0.0198038284421183

$ws.Range("G17").Value = 94.34696966666667
$ws.Range("H17").Value = 283.040909
$ws.Range("I17").Value = 0.524429948940542
$ws.Range("J17").Value = 0.5244299489405421
$ws.Range("M17").Value = 5.912082333333333
$ws.Range("N17").Value = 17.736247
$ws.Range("O17").Value = 0.01926252595751047
$ws.Range("P17").Value = 0.01926252595751047
$ws.Range("Q17").Value = 557.7870525698359
$ws.Range("R17").Value = 5020.083473128522
$ws.Range("S17").Value = 0.01010184550436308
$ws.Range("T17").Value = 0.01010184550436308
